# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - the "BTec_Logo-Orange" picture (appears in the headers)      -> image1.jpg
#   - the Pearson "PearsonLogo.png" picture (appears in the footers) -> image2.png
#
# (Word indexes Headers/Footers as 1=primary, 2=first-page, 3=even-page; this
# walks all three slots of every section so it doesn't matter which slot a
# given logo actually lives in.)

$d = $word.ActiveDocument

function Rename-LogoInStory($story, $newName) {
    if (-not $story.Exists) {
        return
    }
    # Fetching the InlineShape straight off $story.Range can leave a stale
    # handle for footer stories in this host, so re-resolve it through the
    # paragraph that actually contains the drawing.
    $rng = $story.Range
    for ($i = 1; $i -le $rng.Paragraphs.Count; $i++) {
        $para = $rng.Paragraphs($i)
        if ($para.Range.InlineShapes.Count -ge 1) {
            $shape = $para.Range.InlineShapes(1)
            $shape.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($idx = 1; $idx -le 3; $idx++) {
        Rename-LogoInStory $sec.Headers($idx) "image1.jpg"
        Rename-LogoInStory $sec.Footers($idx) "image2.png"
    }
}
